# Auto-generated Excel COM-interop script to apply numeric updates
# to the Odin_Profits workbook (commit: "chore: update Sheets via scheduled runner").
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 5230.5186
$ws.Range("I28").Value = 834.8
$ws.Range("J28").Value = 7816.2354
$ws.Range("K28").Value = 834.8
$ws.Range("L28").Value = 7816.2354
$ws.Range("M28").Value = -349.8
$ws.Range("N28").Value = -8786.2354
# Row 33
$ws.Range("H33").Value = 489.9091
$ws.Range("I33").Value = 489.9091
$ws.Range("K33").Value = 489.9091
$ws.Range("M33").Value = -260.9091
# Row 39
$ws.Range("H39").Value = 99
$ws.Range("I39").Value = 65.333336
$ws.Range("K39").Value = 196.000008
$ws.Range("M39").Value = 99.99999199999999
# Row 40
$ws.Range("H40").Value = 12221.25
$ws.Range("J40").Value = 12221.25
$ws.Range("L40").Value = 12221.25
$ws.Range("N40").Value = -12571.25
# Row 76
$ws.Range("H76").Value = 38468932
$ws.Range("I76").Value = 58831668
$ws.Range("K76").Value = 58831668
$ws.Range("M76").Value = -58831353
# Row 79
$ws.Range("H79").Value = 38468932
$ws.Range("I79").Value = 58831668
$ws.Range("K79").Value = 58831668
$ws.Range("M79").Value = -58830576
# Row 82
$ws.Range("H82").Value = 6666.3335
$ws.Range("I82").Value = 6666.3335
$ws.Range("K82").Value = 19999.0005
$ws.Range("M82").Value = -19593.0005
# Row 85
$ws.Range("H85").Value = 6666.3335
$ws.Range("I85").Value = 6666.3335
$ws.Range("K85").Value = 19999.0005
$ws.Range("M85").Value = -18595.0005
# Row 97
$ws.Range("H97").Value = 1586.6666
$ws.Range("J97").Value = 1586.6666
$ws.Range("L97").Value = 4759.9998
$ws.Range("N97").Value = -5751.9998
# Row 107
$ws.Range("H107").Value = 1437.0555
$ws.Range("I107").Value = 1437.0555
$ws.Range("K107").Value = 1437.0555
$ws.Range("M107").Value = 482.9445000000001
# Row 112
$ws.Range("H112").Value = 2205.6226
$ws.Range("J112").Value = 2307.102
$ws.Range("L112").Value = 6921.306
$ws.Range("N112").Value = -9137.306
# Row 115
$ws.Range("H115").Value = 2995
$ws.Range("J115").Value = 9450
$ws.Range("L115").Value = 28350
$ws.Range("N115").Value = -31484
# Row 118
$ws.Range("H118").Value = 920
$ws.Range("I118").Value = 200
$ws.Range("K118").Value = 600
$ws.Range("M118").Value = 1057
# Row 133
$ws.Range("H133").Value = 74998.336
$ws.Range("J133").Value = 74998.336
$ws.Range("L133").Value = 74998.336
$ws.Range("N133").Value = -85118.336

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 3723.0588
$ws.Range("I2").Value = 1719.4667
$ws.Range("K2").Value = 1719.4667
$ws.Range("M2").Value = -1606.4667
# Row 45
$ws.Range("H45").Value = 1323.3334
$ws.Range("I45").Value = 1298
$ws.Range("K45").Value = 1298
$ws.Range("M45").Value = -921
# Row 110
$ws.Range("H110").Value = 5174.7407
$ws.Range("I110").Value = 2351.9092
$ws.Range("K110").Value = 2351.9092
$ws.Range("M110").Value = -306.9092000000001
# Row 116
$ws.Range("H116").Value = 3723.0588
$ws.Range("I116").Value = 1719.4667
$ws.Range("K116").Value = 1719.4667
$ws.Range("M116").Value = 574.5333000000001

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 3723.0588
$ws.Range("I3").Value = 1719.4667
$ws.Range("K3").Value = 1719.4667
$ws.Range("M3").Value = -1605.4667
# Row 99
$ws.Range("H99").Value = 4976.7437
$ws.Range("I99").Value = 1921.8636
$ws.Range("K99").Value = 1921.8636
$ws.Range("M99").Value = -423.8635999999999
# Row 107
$ws.Range("H107").Value = 2086698.1
$ws.Range("I107").Value = 2859622
$ws.Range("J107").Value = 5749
$ws.Range("K107").Value = 2859622
$ws.Range("L107").Value = 5749
$ws.Range("M107").Value = -2857702
$ws.Range("N107").Value = -9589
# Row 134
$ws.Range("H134").Value = 866662.9
$ws.Range("I134").Value = 896588.3
$ws.Range("K134").Value = 2689764.9
$ws.Range("M134").Value = -2687229.9

$ws = $wb.Worksheets.Item("CRP")
# Row 105
$ws.Range("H105").Value = 40001144
$ws.Range("I105").Value = 55556692
$ws.Range("K105").Value = 55556692
$ws.Range("M105").Value = -55554945
# Row 107
$ws.Range("H107").Value = 1097.3462
$ws.Range("I107").Value = 778.94116
$ws.Range("J107").Value = 1698.7778
$ws.Range("K107").Value = 778.94116
$ws.Range("L107").Value = 1698.7778
$ws.Range("M107").Value = 1141.05884
$ws.Range("N107").Value = -5538.7778

$ws = $wb.Worksheets.Item("CUL")
# Row 121
$ws.Range("H121").Value = 42248.875
$ws.Range("I121").Value = 6500
$ws.Range("K121").Value = 19500
$ws.Range("M121").Value = -18190

$ws = $wb.Worksheets.Item("GSM")
# Row 15
$ws.Range("H15").Value = 9990.666999999999
$ws.Range("J15").Value = 9990.666999999999
$ws.Range("L15").Value = 9990.666999999999
$ws.Range("N15").Value = -10566.667
# Row 81
$ws.Range("H81").Value = 9990.666999999999
$ws.Range("J81").Value = 9990.666999999999
$ws.Range("L81").Value = 9990.666999999999
$ws.Range("N81").Value = -11986.667
# Row 84
$ws.Range("H84").Value = 9990.666999999999
$ws.Range("J84").Value = 9990.666999999999
$ws.Range("L84").Value = 29972.001
$ws.Range("N84").Value = -39956.001

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 981.8333
$ws.Range("I22").Value = 968.6667
$ws.Range("K22").Value = 968.6667
$ws.Range("M22").Value = -673.6667
# Row 27
$ws.Range("H27").Value = 981.8333
$ws.Range("I27").Value = 968.6667
$ws.Range("K27").Value = 968.6667
$ws.Range("M27").Value = -861.6667
# Row 55
$ws.Range("H55").Value = 4229.8423
$ws.Range("J55").Value = 5277.2
$ws.Range("L55").Value = 5277.2
$ws.Range("N55").Value = -5623.2
# Row 61
$ws.Range("H61").Value = 5105.8335
$ws.Range("I61").Value = 3917.074
$ws.Range("K61").Value = 3917.074
$ws.Range("M61").Value = -3715.074
# Row 113
$ws.Range("H113").Value = 5105.8335
$ws.Range("I113").Value = 3917.074
$ws.Range("K113").Value = 3917.074
$ws.Range("M113").Value = -1747.074

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 14837.375
$ws.Range("I62").Value = 12937.25
$ws.Range("J62").Value = 16737.5
$ws.Range("K62").Value = 12937.25
$ws.Range("L62").Value = 16737.5
$ws.Range("M62").Value = -12313.25
$ws.Range("N62").Value = -17985.5
# Row 65
$ws.Range("H65").Value = 14837.375
$ws.Range("I65").Value = 12937.25
$ws.Range("J65").Value = 16737.5
$ws.Range("K65").Value = 64686.25
$ws.Range("L65").Value = 83687.5
$ws.Range("M65").Value = -61566.25
$ws.Range("N65").Value = -89927.5
# Row 107
$ws.Range("H107").Value = 6667503.5
$ws.Range("I107").Value = 11111696
$ws.Range("J107").Value = 1215.8334
$ws.Range("K107").Value = 33335088
$ws.Range("L107").Value = 3647.5002
$ws.Range("M107").Value = -33333168
$ws.Range("N107").Value = -7487.5002
# Row 113
$ws.Range("H113").Value = 9261264
$ws.Range("I113").Value = 13890852
$ws.Range("J113").Value = 2087.3333
$ws.Range("K113").Value = 41672556
$ws.Range("L113").Value = 6261.999899999999
$ws.Range("M113").Value = -41670386
$ws.Range("N113").Value = -10601.9999
# Row 136
$ws.Range("H136").Value = 15165637
$ws.Range("I136").Value = 23821390
$ws.Range("J136").Value = 18070.584
$ws.Range("K136").Value = 71464170
$ws.Range("L136").Value = 54211.75199999999
$ws.Range("M136").Value = -71461620
$ws.Range("N136").Value = -59311.75199999999

